# Re-order the weekly Perejil price records across rows 2-6.
# The new row order (by original source row) is: 4, 6, 3, 2, 5
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original values for the mutable columns (D, J, K, L, M, P) before overwriting.
# Note: use Value2 (not Value) for reads - Value's getter is not usable in this host.
$orig = @{}
foreach ($r in 2..6) {
    $orig[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        P = $ws.Cells.Item($r, 16).Value2
    }
}

# Mapping: new row -> source (old) row
$order = @{ 2 = 4; 3 = 6; 4 = 3; 5 = 2; 6 = 5 }

foreach ($newRow in $order.Keys) {
    $srcRow = $order[$newRow]
    $vals = $orig[$srcRow]
    $ws.Cells.Item($newRow, 4).Value2 = $vals.D
    $ws.Cells.Item($newRow, 10).Value2 = $vals.J
    $ws.Cells.Item($newRow, 11).Value2 = $vals.K
    $ws.Cells.Item($newRow, 12).Value2 = $vals.L
    $ws.Cells.Item($newRow, 13).Value2 = $vals.M
    $ws.Cells.Item($newRow, 16).Value2 = $vals.P
}
